# Actualización automática 2025-10-07 16:30:08
#
# A new sale of 6777.81 (October) is recorded for client
# "PROVEEDORA PARA METALMECANICA E INDUSTRIAS PROMETIN CIA LTDA"
# (advisor ALMEIDA CUATIN JHONATHANN CARLOS) in the PORCELANATO group.
# This ripples through the three report sheets that were previously
# maintained by hand / by a refresh macro.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" — per-client sales broken out by product
# group. Column M is PORCELANATO. Row 26 is the client above.
# Row 36 is the "N de 34" non-zero-count summary row.
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M26").Value = 6777.81
$wsGrupo.Range("M36").Value = "3 de 34"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" — per-client sales broken out by month.
# Column F is octubre (October). Row 26 is the same client; row 36
# is the column-total row.
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F26").Value = 6777.81
$wsMensual.Range("F36").Value = 8373

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" — budget vs. sale per product group.
# Row 12 is PORCELANATO, row 14 is the TOTAL row. D = VENTA,
# E = POR CUMPLIR (C - D), F = CUMPLIMIENTO (D / C).
# ---------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D widened from 12 to 13 characters. The engine's stored
# XML "width" unit runs 0.8333333333333333 (=5/6) wider than the
# ColumnWidth value you assign, so back that constant out here to
# land exactly on width="13" in the saved file.
$wsCumplimiento.Range("D1").ColumnWidth = 13 - 0.8333333333333333

$wsCumplimiento.Range("D12").Value = 7617.31
$wsCumplimiento.Range("E12").Value = 14083.96
$wsCumplimiento.Range("F12").Value = 0.3510075677598593

$wsCumplimiento.Range("D14").Value = 8373
$wsCumplimiento.Range("E14").Value = 28212.56723718183
$wsCumplimiento.Range("F14").Value = 0.2288607402399528
